$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 465.66666
$ws.Range("J2").Value = 831.3333
$ws.Range("L2").Value = 831.3333
$ws.Range("N2").Value = -1057.3333

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H51").Value = 4372.5713
$ws.Range("I51").Value = 5040
$ws.Range("K51").Value = 5040
$ws.Range("M51").Value = -4556

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 8261.471
$ws.Range("J69").Value = 9449.846
$ws.Range("L69").Value = 28349.538
$ws.Range("N69").Value = -30097.538

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H72").Value = 8261.471
$ws.Range("J72").Value = 9449.846
$ws.Range("L72").Value = 85048.614
$ws.Range("N72").Value = -93784.614

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1202
$ws.Range("I112").Value = 0
$ws.Range("J112").Value = 1202
$ws.Range("K112").Value = 0
$ws.Range("L112").Value = 3606
$ws.Range("N112").Value = -5822
$ws.Range("M112").ClearContents()

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H113").Value = 2337.647
$ws.Range("I113").Value = 2309.3
$ws.Range("J113").Value = 2378.1428
$ws.Range("K113").Value = 2309.3
$ws.Range("L113").Value = 2378.1428
$ws.Range("M113").Value = 944.6999999999998
$ws.Range("N113").Value = -8886.1428

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H115").Value = 805.6
$ws.Range("I115").Value = 805.6
$ws.Range("K115").Value = 2416.8
$ws.Range("M115").Value = -849.8000000000002

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 3317.2144
$ws.Range("I132").Value = 3205.1667
$ws.Range("K132").Value = 9615.500100000001
$ws.Range("M132").Value = -7085.500100000001

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H141").Value = 1537.8
$ws.Range("I141").Value = 1586.5555
$ws.Range("K141").Value = 4759.666499999999
$ws.Range("M141").Value = 420.3335000000006

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4670.387
$ws.Range("I32").Value = 3490.9167
$ws.Range("J32").Value = 8714.286
$ws.Range("K32").Value = 3490.9167
$ws.Range("L32").Value = 8714.286
$ws.Range("M32").Value = -3203.9167
$ws.Range("N32").Value = -9288.286

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2607.04
$ws.Range("I45").Value = 2940.5881
$ws.Range("K45").Value = 2940.5881
$ws.Range("M45").Value = -2563.5881

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4474.7827
$ws.Range("I61").Value = 1693.4615
$ws.Range("J61").Value = 8090.5
$ws.Range("K61").Value = 1693.4615
$ws.Range("L61").Value = 8090.5
$ws.Range("M61").Value = -1481.4615
$ws.Range("N61").Value = -8514.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H129").Value = 70000
$ws.Range("J129").Value = 70000
$ws.Range("L129").Value = 70000
$ws.Range("N129").Value = -80000

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 6080.8945
$ws.Range("I132").Value = 4410.6924
$ws.Range("J132").Value = 9699.666999999999
$ws.Range("K132").Value = 13232.0772
$ws.Range("L132").Value = 29099.001
$ws.Range("M132").Value = -10702.0772
$ws.Range("N132").Value = -34159.001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4474.7827
$ws.Range("I136").Value = 1693.4615
$ws.Range("J136").Value = 8090.5
$ws.Range("K136").Value = 5080.3845
$ws.Range("L136").Value = 24271.5
$ws.Range("M136").Value = -2530.3845
$ws.Range("N136").Value = -29371.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3441.9512
$ws.Range("I134").Value = 1548.1111
$ws.Range("J134").Value = 7094.357
$ws.Range("K134").Value = 4644.3333
$ws.Range("L134").Value = 21283.071
$ws.Range("M134").Value = -2109.3333
$ws.Range("N134").Value = -26353.071

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 3422.6667
$ws.Range("I16").Value = 3400.7144
$ws.Range("J16").Value = 3499.5
$ws.Range("K16").Value = 3400.7144
$ws.Range("L16").Value = 3499.5
$ws.Range("M16").Value = -3113.7144
$ws.Range("N16").Value = -4073.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 4883.5
$ws.Range("I99").Value = 4854.9033
$ws.Range("J99").Value = 4930.1577
$ws.Range("K99").Value = 4854.9033
$ws.Range("L99").Value = 4930.1577
$ws.Range("M99").Value = -3356.9033
$ws.Range("N99").Value = -7926.1577

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H107").Value = 2034.9286
$ws.Range("I107").Value = 1999.1538
$ws.Range("K107").Value = 1999.1538
$ws.Range("M107").Value = -79.15380000000005

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H113").Value = 3422.6667
$ws.Range("I113").Value = 3400.7144
$ws.Range("J113").Value = 3499.5
$ws.Range("K113").Value = 3400.7144
$ws.Range("L113").Value = 3499.5
$ws.Range("M113").Value = -1230.7144
$ws.Range("N113").Value = -7839.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 4883.5
$ws.Range("I126").Value = 4854.9033
$ws.Range("J126").Value = 4930.1577
$ws.Range("K126").Value = 14564.7099
$ws.Range("L126").Value = 14790.4731
$ws.Range("M126").Value = -12094.7099
$ws.Range("N126").Value = -19730.4731

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H131").Value = 0
$ws.Range("J131").Value = 0
$ws.Range("L131").Value = 0
$ws.Range("N131").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2708.5
$ws.Range("I132").Value = 2122.7144
$ws.Range("J132").Value = 3528.6
$ws.Range("K132").Value = 6368.1432
$ws.Range("L132").Value = 10585.8
$ws.Range("M132").Value = -3838.1432
$ws.Range("N132").Value = -15645.8

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H133").Value = 61998.168
$ws.Range("J133").Value = 61998.168
$ws.Range("L133").Value = 61998.168
$ws.Range("N133").Value = -67058.16800000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H47").Value = 722.25
$ws.Range("I47").Value = 863
$ws.Range("K47").Value = 2589
$ws.Range("M47").Value = -2158

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H60").Value = 2596.9333
$ws.Range("I60").Value = 2055.5
$ws.Range("J60").Value = 3679.8
$ws.Range("K60").Value = 6166.5
$ws.Range("L60").Value = 11039.4
$ws.Range("M60").Value = -5915.5
$ws.Range("N60").Value = -11541.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H102").Value = 7006
$ws.Range("I102").Value = 1737.5
$ws.Range("J102").Value = 10016.571
$ws.Range("K102").Value = 5212.5
$ws.Range("L102").Value = 30049.713
$ws.Range("M102").Value = -2778.5
$ws.Range("N102").Value = -34917.713

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H108").Value = 8755.833000000001
$ws.Range("I108").Value = 1575
$ws.Range("J108").Value = 15936.667
$ws.Range("K108").Value = 4725
$ws.Range("L108").Value = 47810.001
$ws.Range("M108").Value = -1845
$ws.Range("N108").Value = -53570.001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H112").Value = 750
$ws.Range("I112").Value = 750
$ws.Range("J112").Value = 0
$ws.Range("K112").Value = 2250
$ws.Range("L112").Value = 0
$ws.Range("M112").Value = -1142
$ws.Range("N112").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H132").Value = 1675
$ws.Range("I132").Value = 1335.2
$ws.Range("J132").Value = 3374
$ws.Range("K132").Value = 12016.8
$ws.Range("L132").Value = 30366
$ws.Range("M132").Value = -9486.800000000001
$ws.Range("N132").Value = -35426

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 78999
$ws.Range("J42").Value = 78999
$ws.Range("L42").Value = 78999
$ws.Range("N42").Value = -79969

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 0
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 0
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 0
$ws.Range("M70").ClearContents()
$ws.Range("N70").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H73").Value = 0
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 0
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 0
$ws.Range("M73").ClearContents()
$ws.Range("N73").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 10951.25
$ws.Range("I80").Value = 7999
$ws.Range("J80").Value = 11935.333
$ws.Range("K80").Value = 7999
$ws.Range("L80").Value = 11935.333
$ws.Range("M80").Value = -7001
$ws.Range("N80").Value = -13931.333

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H83").Value = 10951.25
$ws.Range("I83").Value = 7999
$ws.Range("J83").Value = 11935.333
$ws.Range("K83").Value = 39995
$ws.Range("L83").Value = 59676.665
$ws.Range("M83").Value = -35003
$ws.Range("N83").Value = -69660.66500000001

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 1383.875
$ws.Range("I107").Value = 1040.1818
$ws.Range("K107").Value = 1040.1818
$ws.Range("M107").Value = 879.8181999999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H115").Value = 78999
$ws.Range("J115").Value = 78999
$ws.Range("L115").Value = 78999
$ws.Range("N115").Value = -81349

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 4263.472
$ws.Range("I122").Value = 2660.8696
$ws.Range("J122").Value = 7098.846
$ws.Range("K122").Value = 7982.6088
$ws.Range("L122").Value = 21296.538
$ws.Range("M122").Value = -5532.6088
$ws.Range("N122").Value = -26196.538

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 2866.182
$ws.Range("I126").Value = 2652.8
$ws.Range("J126").Value = 5000
$ws.Range("K126").Value = 7958.400000000001
$ws.Range("L126").Value = 15000
$ws.Range("M126").Value = -5488.400000000001
$ws.Range("N126").Value = -19940

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 1887.35
$ws.Range("I132").Value = 1658.7778
$ws.Range("K132").Value = 4976.3334
$ws.Range("M132").Value = -2446.3334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H24").Value = 20000
$ws.Range("J24").Value = 20000
$ws.Range("L24").Value = 20000
$ws.Range("N24").Value = -20686

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H58").Value = 6298.5
$ws.Range("I58").Value = 7800
$ws.Range("K58").Value = 7800
$ws.Range("M58").Value = -7540

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2652.439
$ws.Range("I122").Value = 2395.7144
$ws.Range("K122").Value = 7187.1432
$ws.Range("M122").Value = -4737.1432

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H20").Value = 0
$ws.Range("I20").Value = 0
$ws.Range("K20").Value = 0
$ws.Range("M20").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4466.0835
$ws.Range("J96").Value = 4499.5
$ws.Range("L96").Value = 4499.5
$ws.Range("N96").Value = -7245.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4106.659
$ws.Range("I132").Value = 2378.375
$ws.Range("K132").Value = 7135.125
$ws.Range("M132").Value = -4605.125
